# Sync automático del tracker - 2025-08-17 21:40:10 UTC
#
# 1) Predictions sheet: mark the 2025-08-17 Real Sociedad II vs Zaragoza
#    pick as resolved (Status Pending -> Completed).
# 2) Results sheet: append the corresponding settled-result row (row 30)
#    now that the match has finished.

$wb = $excel.ActiveWorkbook

$predictions = $wb.Worksheets.Item("Predictions")
$predictions.Range("I30").Value = "Completed"

$results = $wb.Worksheets.Item("Results")

# Dates in this tracker are stored as plain text (e.g. "2025-08-17"),
# not native Excel date serials. Force text format first so the engine
# doesn't auto-coerce the string into a date value, then drop back to an
# unstyled "Normal" cell so no stray number-format style sticks around.
$results.Range("A30").NumberFormat = "@"
$results.Cells.Item(30, 1).Value = "2025-08-17"
$results.Range("A30").Style = "Normal"

$results.Cells.Item(30, 2).Value = "Segunda División"
$results.Cells.Item(30, 3).Value = "real sociedad ii"
$results.Cells.Item(30, 4).Value = "zaragoza"
$results.Cells.Item(30, 5).Value = "Home Win"
$results.Cells.Item(30, 6).Value = "Away Win"
$results.Cells.Item(30, 7).Value = $false
$results.Cells.Item(30, 8).Value = -1
$results.Cells.Item(30, 9).Value = -100

$results.Range("J30").NumberFormat = "@"
$results.Cells.Item(30, 10).Value = "2025-08-17"
$results.Range("J30").Style = "Normal"
